$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New employee records appended to Sheet1 (rows 14-16), mirroring the
# existing "IRCMP" block (row 13) but for a new employee "Doe, John".
$newRows = @(
    @{ Row = 14; B = 10323194; K = 3 },
    @{ Row = 15; B = 10473193; K = 4 },
    @{ Row = 16; B = 10411336; K = 5 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value2 = "IRCMP"
    $ws.Range("B$row").Value2 = $r.B
    $ws.Range("C$row").Value2 = "Doe, John"
    $ws.Range("D$row").Value2 = "A"
    $ws.Range("E$row").Value2 = "Short Work Break"
    $ws.Range("F$row").Value2 = "9AC"
    $ws.Range("G$row").Value2 = 45931
    $ws.Range("H$row").Value2 = 45912
    $ws.Range("I$row").Value2 = "F"
    $ws.Range("J$row").Value2 = "F"
    $ws.Range("K$row").Value2 = $r.K
    $ws.Range("L$row").Value2 = "P"

    # M must stay the text value "1029075745931" (same text already used by
    # row 13), not get auto-coerced to a number. Copying the value from M13
    # (values only) keeps it as text without introducing any new style.
    $ws.Range("M13").Copy() | Out-Null
    $ws.Range("M$row").PasteSpecial(-4163) | Out-Null

    # Copy formatting (not values) from row 13, which already carries the
    # correct styles for every column except B and C.
    $ws.Range("A13:M13").Copy() | Out-Null
    $ws.Range("A$row`:M$row").PasteSpecial(-4122) | Out-Null

    # Column C needs the regular (non-bold) style used elsewhere in the
    # sheet, e.g. C5, rather than row 13's bold style.
    $ws.Range("C5").Copy() | Out-Null
    $ws.Range("C$row").PasteSpecial(-4122) | Out-Null

    # Column B keeps the default (unstyled) formatting, matching B2.
    $ws.Range("B$row").ClearFormats() | Out-Null
}

$excel.CutCopyMode = $false

$ws.Range("D1:L1048576").Select() | Out-Null
$ws.Application.ActiveCell = $ws.Range("D1")
